# Apply updated market-price figures to the leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3132.1333
$ws.Range("I32").Value = 4048.3333
$ws.Range("J32").Value = 2521.3333
$ws.Range("K32").Value = 4048.3333
$ws.Range("L32").Value = 2521.3333
$ws.Range("M32").Value = -3722.3333
$ws.Range("N32").Value = -3173.3333

$ws.Range("H111").Value = 2646.75
$ws.Range("I111").Value = 2435.818
$ws.Range("J111").Value = 3110.8
$ws.Range("K111").Value = 7307.454000000001
$ws.Range("L111").Value = 9332.400000000001
$ws.Range("M111").Value = -4240.454000000001
$ws.Range("N111").Value = -15466.4

$ws.Range("H132").Value = 4693.0835
$ws.Range("I132").Value = 4693.0835
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14079.2505
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11549.2505
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 866882.5
$ws.Range("I2").Value = 1700.7894
$ws.Range("J2").Value = 1962779.2
$ws.Range("K2").Value = 1700.7894
$ws.Range("L2").Value = 1962779.2
$ws.Range("M2").Value = -1587.7894
$ws.Range("N2").Value = -1963005.2

$ws.Range("H61").Value = 8549400
$ws.Range("I61").Value = 9261600
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 9261600
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -9261388
$ws.Range("N61").Value = -3424

$ws.Range("H116").Value = 866882.5
$ws.Range("I116").Value = 1700.7894
$ws.Range("J116").Value = 1962779.2
$ws.Range("K116").Value = 1700.7894
$ws.Range("L116").Value = 1962779.2
$ws.Range("M116").Value = 593.2106000000001
$ws.Range("N116").Value = -1967367.2

$ws.Range("H132").Value = 1778331.2
$ws.Range("I132").Value = 2467338
$ws.Range("K132").Value = 7402014
$ws.Range("M132").Value = -7399484

$ws.Range("H136").Value = 8549400
$ws.Range("I136").Value = 9261600
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 27784800
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -27782250
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 866882.5
$ws.Range("I3").Value = 1700.7894
$ws.Range("J3").Value = 1962779.2
$ws.Range("K3").Value = 1700.7894
$ws.Range("L3").Value = 1962779.2
$ws.Range("M3").Value = -1586.7894
$ws.Range("N3").Value = -1963007.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 16718.3
$ws.Range("J59").Value = 16798.111
$ws.Range("L59").Value = 16798.111
$ws.Range("N59").Value = -19088.111

$ws.Range("H107").Value = 3473356.2
$ws.Range("I107").Value = 5953542
$ws.Range("J107").Value = 1096
$ws.Range("K107").Value = 5953542
$ws.Range("L107").Value = 1096
$ws.Range("M107").Value = -5951622
$ws.Range("N107").Value = -4936

$ws.Range("H141").Value = 29299.7
$ws.Range("J141").Value = 29299.7
$ws.Range("L141").Value = 29299.7
$ws.Range("N141").Value = -39659.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 9402068
$ws.Range("I122").Value = 17544422
$ws.Range("J122").Value = 1259714.1
$ws.Range("K122").Value = 157899798
$ws.Range("L122").Value = 11337426.9
$ws.Range("M122").Value = -157897348
$ws.Range("N122").Value = -11342326.9

$ws.Range("H131").Value = 5619.341
$ws.Range("I131").Value = 4163.75
$ws.Range("J131").Value = 6451.107
$ws.Range("K131").Value = 12491.25
$ws.Range("L131").Value = 19353.321
$ws.Range("M131").Value = -7451.25
$ws.Range("N131").Value = -29433.321

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 24499.5
$ws.Range("J34").Value = 24499.5
$ws.Range("L34").Value = 24499.5
$ws.Range("N34").Value = -25035.5

$ws.Range("H76").Value = 24499.5
$ws.Range("J76").Value = 24499.5
$ws.Range("L76").Value = 24499.5
$ws.Range("N76").Value = -25129.5

$ws.Range("H79").Value = 24499.5
$ws.Range("J79").Value = 24499.5
$ws.Range("L79").Value = 24499.5
$ws.Range("N79").Value = -26683.5

$ws.Range("H113").Value = 1635.0834
$ws.Range("I113").Value = 1349
$ws.Range("J113").Value = 2493.3333
$ws.Range("K113").Value = 1349
$ws.Range("L113").Value = 2493.3333
$ws.Range("M113").Value = 821
$ws.Range("N113").Value = -6833.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H20").Value = 3333.3333
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5452

$ws.Range("H61").Value = 1476.5714
$ws.Range("I61").Value = 1476.5714
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1476.5714
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1274.5714
$ws.Range("N61").ClearContents()

$ws.Range("H93").Value = 1148.85
$ws.Range("I93").Value = 754.05554
$ws.Range("J93").Value = 4702
$ws.Range("K93").Value = 754.05554
$ws.Range("L93").Value = 4702
$ws.Range("M93").Value = 493.94446
$ws.Range("N93").Value = -7198

$ws.Range("H100").Value = 2603.348
$ws.Range("I100").Value = 1832.5
$ws.Range("J100").Value = 3444.2727
$ws.Range("K100").Value = 1832.5
$ws.Range("L100").Value = 3444.2727
$ws.Range("M100").Value = -1291.5
$ws.Range("N100").Value = -4526.2727

$ws.Range("H103").Value = 13620
$ws.Range("J103").Value = 13620
$ws.Range("L103").Value = 13620
$ws.Range("N103").Value = -15964

$ws.Range("H113").Value = 1476.5714
$ws.Range("I113").Value = 1476.5714
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1476.5714
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 693.4286
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 551.8570999999999
$ws.Range("I107").Value = 395.0909
$ws.Range("J107").Value = 724.3
$ws.Range("K107").Value = 1185.2727
$ws.Range("L107").Value = 2172.9
$ws.Range("M107").Value = 734.7273
$ws.Range("N107").Value = -6012.9

$ws.Range("H113").Value = 392.95
$ws.Range("I113").Value = 393.07144
$ws.Range("J113").Value = 392.66666
$ws.Range("K113").Value = 1179.21432
$ws.Range("L113").Value = 1177.99998
$ws.Range("M113").Value = 990.78568
$ws.Range("N113").Value = -5517.999980000001
